$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.032160878831703751
$ws.Range("B1").Value = -0.032223395104730081
$ws.Range("A2").Value = -0.013920030670878566
$ws.Range("B2").Value = -0.021384908622313414
$ws.Range("A3").Value = -0.069483840199985783
$ws.Range("B3").Value = -0.069499642919175958
